$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.925.80'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.662.43'
$ws.Range("E3").Value = '  -1.92%  '
$ws.Range("D4").Value = "'0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.65%  '
$ws.Range("D5").Value = "'317.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").Value = "'0.9966"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").Value = "'0.3638"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.12%  '
$ws.Range("D8").Value = "'47.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.30%  '
$ws.Range("D9").Value = "'0.3279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.80%  '
$ws.Range("D10").Value = "'1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07083"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.61%  '
$ws.Range("D12").Value = "'0.9955"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = "'6.063"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.58%  '
$ws.Range("D14").Value = "'19.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.61%  '
$ws.Range("D15").Value = '1.666.38'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = "'6.636"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("D17").Value = "'0.00001053"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.74%  '
$ws.Range("D18").Value = "'0.06626"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = "'0.9963"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("D20").Value = "'79.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.14%  '
$ws.Range("D21").Value = "'5.936"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.63%  '
$ws.Range("D22").Value = "'15.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.51%  '
$ws.Range("D23").Value = "'12.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").Value = '24.883.63'
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").Value = "'2.434"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").Value = "'2.415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -12.24%  '
$ws.Range("D27").Value = "'148.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").Value = "'18.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.94%  '
$ws.Range("D29").Value = "'1.240"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = '1.849.63'
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("D31").Value = "'126.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("D32").Value = "'4.097"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.95%  '
$ws.Range("D33").Value = "'5.863"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -12.34%  '
$ws.Range("D34").Value = "'0.08454"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").Value = "'1.679"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.00%  '
$ws.Range("D36").Value = "'12.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.04%  '
$ws.Range("D37").Value = "'1.295"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.22%  '
$ws.Range("D38").Value = "'5.235"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.96%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.02254"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.71%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.06050"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.00%  '
$ws.Range("D41").Value = "'0.2077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.52%  '
$ws.Range("D42").Value = "'8.281"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.22%  '
$ws.Range("D43").Value = "'0.9954"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").Value = "'0.5958"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.34%  '
$ws.Range("D45").Value = "'3.827"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").Value = "'12.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.96%  '
$ws.Range("D47").Value = "'0.5667"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.09%  '
$ws.Range("D48").Value = "'125.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = "'1.961"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.48%  '
$ws.Range("D50").Value = "'0.07019"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").Value = "'1.199"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.98%  '
